$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.246.18"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.541.45"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "3.541.36"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.12"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "4.147.14"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.49"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "3.547.84"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "65.338.71"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.79"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.73"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.581"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").Value = "3.686.62"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.66"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.89"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  +16.45%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.30"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.48"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("D33").Value = "3.550.68"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.147"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.31"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "168.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.89"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.827"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.90"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Value = "2.395.39"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0268"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.49%  "
